# Change reference.docx to use more normal block quotes: indent the
# "Block Text" paragraph style on the left and right by 480 twips (24pt),
# matching the commit's updated <w:ind w:firstLine="0" w:left="480"
# w:right="480"/> on w:style[@w:styleId='BlockText']/w:pPr.

$d = $word.ActiveDocument

$blockText = $d.Styles("Block Text")

$blockText.ParagraphFormat.LeftIndent = 24   # points -> 480 twips
$blockText.ParagraphFormat.RightIndent = 24  # points -> 480 twips
